$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 559
$ws.Range("I6").Value = 316.66666
$ws.Range("J6").Value = 1649.5
$ws.Range("K6").Value = 949.9999799999999
$ws.Range("L6").Value = 4948.5
$ws.Range("M6").Value = -837.9999799999999
$ws.Range("N6").Value = -5172.5
$ws.Range("H34").Value = 3418.8572
$ws.Range("I34").Value = 655.3333
$ws.Range("J34").Value = 20000
$ws.Range("K34").Value = 655.3333
$ws.Range("L34").Value = 20000
$ws.Range("M34").Value = -452.3333
$ws.Range("N34").Value = -20406
$ws.Range("H36").Value = 3418.8572
$ws.Range("I36").Value = 655.3333
$ws.Range("J36").Value = 20000
$ws.Range("K36").Value = 655.3333
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = 59.66669999999999
$ws.Range("N36").Value = -21430
$ws.Range("H64").Value = 3358.383
$ws.Range("I64").Value = 2998.6572
$ws.Range("J64").Value = 4407.5835
$ws.Range("K64").Value = 2998.6572
$ws.Range("L64").Value = 4407.5835
$ws.Range("M64").Value = -2750.6572
$ws.Range("N64").Value = -4903.5835
$ws.Range("H67").Value = 3358.383
$ws.Range("I67").Value = 2998.6572
$ws.Range("J67").Value = 4407.5835
$ws.Range("K67").Value = 2998.6572
$ws.Range("L67").Value = 4407.5835
$ws.Range("M67").Value = -2140.6572
$ws.Range("N67").Value = -6123.5835
$ws.Range("H74").Value = 4429
$ws.Range("I74").Value = 5003
$ws.Range("J74").Value = 4333.3335
$ws.Range("K74").Value = 5003
$ws.Range("L74").Value = 4333.3335
$ws.Range("M74").Value = -4067
$ws.Range("N74").Value = -6205.3335
$ws.Range("H77").Value = 4429
$ws.Range("I77").Value = 5003
$ws.Range("J77").Value = 4333.3335
$ws.Range("K77").Value = 25015
$ws.Range("L77").Value = 21666.6675
$ws.Range("M77").Value = -20335
$ws.Range("N77").Value = -31026.6675
$ws.Range("H93").Value = 200000
$ws.Range("J93").Value = 200000
$ws.Range("L93").Value = 200000
$ws.Range("N93").Value = -204992
$ws.Range("H112").Value = 1141.3334
$ws.Range("J112").Value = 1144
$ws.Range("L112").Value = 3432
$ws.Range("N112").Value = -5648
$ws.Range("H123").Value = 10000
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 21394.922
$ws.Range("I132").Value = 26851.125
$ws.Range("J132").Value = 1554.1818
$ws.Range("K132").Value = 80553.375
$ws.Range("L132").Value = 4662.5454
$ws.Range("M132").Value = -78023.375
$ws.Range("N132").Value = -9722.545399999999
$ws.Range("H137").Value = 1152.6545
$ws.Range("I137").Value = 999.6429000000001
$ws.Range("J137").Value = 1647
$ws.Range("K137").Value = 2998.9287
$ws.Range("L137").Value = 4941
$ws.Range("M137").Value = -448.9287000000004
$ws.Range("N137").Value = -10041
$ws.Range("H138").Value = 3393.8572
$ws.Range("J138").Value = 4899.718
$ws.Range("L138").Value = 14699.154
$ws.Range("N138").Value = -24979.154

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 15913.8
$ws.Range("I6").Value = 19725.75
$ws.Range("J6").Value = 666
$ws.Range("K6").Value = 19725.75
$ws.Range("L6").Value = 666
$ws.Range("M6").Value = -19552.75
$ws.Range("N6").Value = -1012
$ws.Range("H45").Value = 2344.8572
$ws.Range("I45").Value = 2841.3333
$ws.Range("K45").Value = 2841.3333
$ws.Range("M45").Value = -2464.3333
$ws.Range("H101").Value = 49301
$ws.Range("J101").Value = 49301
$ws.Range("L101").Value = 49301
$ws.Range("N101").Value = -55791
$ws.Range("H109").Value = 24800
$ws.Range("J109").Value = 24800
$ws.Range("L109").Value = 24800
$ws.Range("N109").Value = -27574

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 29500
$ws.Range("J64").Value = 29500
$ws.Range("L64").Value = 29500
$ws.Range("N64").Value = -29996
$ws.Range("H67").Value = 29500
$ws.Range("J67").Value = 29500
$ws.Range("L67").Value = 29500
$ws.Range("N67").Value = -31216
$ws.Range("H99").Value = 13632.8
$ws.Range("I99").Value = 3000
$ws.Range("J99").Value = 18189.715
$ws.Range("K99").Value = 3000
$ws.Range("L99").Value = 18189.715
$ws.Range("M99").Value = -1502
$ws.Range("N99").Value = -21185.715
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H126").Value = 13632.8
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 18189.715
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 54569.145
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -59509.145
$ws.Range("H132").Value = 2021.9166
$ws.Range("I132").Value = 1293.0667
$ws.Range("K132").Value = 3879.2001
$ws.Range("M132").Value = -1349.2001

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5486.5625
$ws.Range("I131").Value = 734
$ws.Range("K131").Value = 2202
$ws.Range("M131").Value = 2838

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 5464444.5
$ws.Range("I12").Value = 5272500
$ws.Range("J12").Value = 7000000
$ws.Range("K12").Value = 5272500
$ws.Range("L12").Value = 7000000
$ws.Range("M12").Value = -5272360
$ws.Range("N12").Value = -7000280
$ws.Range("H80").Value = 1998.1428
$ws.Range("I80").Value = 1996.875
$ws.Range("J80").Value = 1999.8334
$ws.Range("K80").Value = 1996.875
$ws.Range("L80").Value = 1999.8334
$ws.Range("M80").Value = -998.875
$ws.Range("N80").Value = -3995.8334
$ws.Range("H83").Value = 1998.1428
$ws.Range("I83").Value = 1996.875
$ws.Range("J83").Value = 1999.8334
$ws.Range("K83").Value = 9984.375
$ws.Range("L83").Value = 9999.166999999999
$ws.Range("M83").Value = -4992.375
$ws.Range("N83").Value = -19983.167

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 69646.664
$ws.Range("I7").Value = 102330
$ws.Range("J7").Value = 4280
$ws.Range("K7").Value = 102330
$ws.Range("L7").Value = 4280
$ws.Range("M7").Value = -102218
$ws.Range("N7").Value = -4504
$ws.Range("H104").Value = 11600
$ws.Range("J104").Value = 11600
$ws.Range("L104").Value = 11600
$ws.Range("N104").Value = -18588
$ws.Range("H126").Value = 69646.664
$ws.Range("I126").Value = 102330
$ws.Range("J126").Value = 4280
$ws.Range("K126").Value = 306990
$ws.Range("L126").Value = 12840
$ws.Range("M126").Value = -304520
$ws.Range("N126").Value = -17780

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 25333.666
$ws.Range("J101").Value = 25333.666
$ws.Range("L101").Value = 25333.666
$ws.Range("N101").Value = -31823.666
$ws.Range("H126").Value = 50991.355
$ws.Range("I126").Value = 59275.125
$ws.Range("J126").Value = 1288.75
$ws.Range("K126").Value = 177825.375
$ws.Range("L126").Value = 3866.25
$ws.Range("M126").Value = -175355.375
$ws.Range("N126").Value = -8806.25
$ws.Range("H132").Value = 1449.9656
$ws.Range("I132").Value = 1292.7778
$ws.Range("J132").Value = 1994.0769
$ws.Range("K132").Value = 3878.3334
$ws.Range("L132").Value = 5982.2307
$ws.Range("M132").Value = -1348.3334
$ws.Range("N132").Value = -11042.2307

Write-Output "Applied all Carbuncle_Profits cell updates"